# Excel COM-interop script replicating the "cartoons" workbook update:
#  - Mark EndYear ("E") column with right alignment, and fill ongoing
#    shows' EndYear with "Present"
#  - Add header alignment (left) on E1
#  - Append 6 new cartoon rows (18-23) with full data
#  - Add Hyperlink-styled Image cells for the new rows, with real
#    hyperlinks for 3 of them
#  - Resize column A, move the active selection, and tidy up

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlRight = -4152
$xlLeft  = -4131

# ---------------------------------------------------------------------
# 1. Existing data rows (2-17): style the EndYear column, and fill in
#    "Present" for shows that are still ongoing (previously blank E).
# ---------------------------------------------------------------------
$presentRows = @(5,6,8,9,10,13,14)
for ($r = 2; $r -le 17; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    if ($presentRows -contains $r) {
        $cell.Value = "Present"
    }
    $cell.HorizontalAlignment = $xlRight
}

# Header cell E1 keeps its bold font but becomes left-aligned.
$ws.Cells.Item(1, 5).HorizontalAlignment = $xlLeft

# ---------------------------------------------------------------------
# 2. Append the 6 new rows of show data.
# ---------------------------------------------------------------------
function Add-CartoonRow {
    param($row, $title, $creator, $desc, $startYear, $endYear, $genre,
          $seasons, $episodes, $image, $network, $hyperlink)

    $ws.Cells.Item($row, 1).Value = $title
    $ws.Cells.Item($row, 2).Value = $creator
    $ws.Cells.Item($row, 3).Value = $desc
    $ws.Cells.Item($row, 4).Value = $startYear

    $eCell = $ws.Cells.Item($row, 5)
    if ($endYear -eq "Present") {
        $eCell.Value = "Present"
    } else {
        $eCell.Value = $endYear
    }
    $eCell.HorizontalAlignment = $xlRight

    $ws.Cells.Item($row, 6).Value = $genre
    $ws.Cells.Item($row, 7).Value = $seasons
    $ws.Cells.Item($row, 8).Value = $episodes

    $iCell = $ws.Cells.Item($row, 9)
    $iCell.Value = $image
    if ($hyperlink) {
        $ws.Hyperlinks.Add($iCell, $image) | Out-Null
    }
    $iCell.Style = "Hyperlink"

    $ws.Cells.Item($row, 10).Value = $network
}

Add-CartoonRow 18 "Randy Cunningham: 9th Grade Ninja" "Jed Elinoff, Scott Thomas" `
    "An ordinary ninth grader is chosen to become The Ninja. He is tasked with protecting the town from evil as well as balancing his school life with friends and homework." `
    2012 2015 "Action, Comedy" 2 100 `
    "https://m.media-amazon.com/images/M/MV5BNTIzMzliOWItNDNhYi00N2YyLWI1NTQtMjMxMTc1ZjhlMDk4XkEyXkFqcGdeQXVyMzM4NjcxOTc@._V1_SY1000_SX642_AL_.jpg" `
    "Disney XD" $false

Add-CartoonRow 19 "Penn Zero: Part Time Hero" "Jared Bush, Sam Levine" `
    "Penn Zero is not your average kid - every day he's zapped into another dimension with his friends to save the world." `
    2014 2017 "Action, Adventure, Science Fantasy" 2 61 `
    "https://m.media-amazon.com/images/M/MV5BN2VmMDgyODAtYWI3Ni00NWZiLTkxODktZTJlYTc0MjA1M2UxXkEyXkFqcGdeQXVyMzM4NjcxOTc@._V1_SY1000_CR0,0,703,1000_AL_.jpg" `
    "Disney XD" $false

Add-CartoonRow 20 "Milo Murphy's Law" 'Dan Povenmire, Jeff "Swampy" Marsh' `
    "An animated comedy adventure series that follows 13-year-old Milo Murphy, the fictional great-great-great-great grandson of the Murphy's Law namesake." `
    2016 "Present" "Comedy" 2 40 `
    "https://m.media-amazon.com/images/M/MV5BMjQxMDY0NjY1MV5BMl5BanBnXkFtZTgwNzQwNDc4OTE@._V1_SY1000_SX690_AL_.jpg" `
    "Disney XD" $false

Add-CartoonRow 21 "The Adventures of Jimmy Neutron: Boy Genius" "John A. Davis, Keith Alcorn" `
    "A young boy, who happens to be a genius, lives in a small town with his family and friends and often gets into crazy adventures with them involving the things he invents." `
    2002 2006 "Adventure, Science Fiction" 3 81 `
    "https://m.media-amazon.com/images/M/MV5BMWRlNTRkM2ItNDkwMC00ZjNmLWI2ZDQtNWI0MTllMGU5OTVjXkEyXkFqcGdeQXVyNTUyMzE4Mzg@._V1_SY1000_CR0,0,750,1000_AL_.jpg" `
    "Nickelodeon" $true

Add-CartoonRow 22 "Camp Lazlo" "Joe Murray" `
    "Bean Scout Lazlo, a fun-loving, free-spirited monkey, and his two bunkmates Raj and Clam, wreak havoc on a very structured summer camp." `
    2005 2008 "Comedy, Slapstick" 5 120 `
    "https://m.media-amazon.com/images/M/MV5BNDQwYzI3YzQtZTQyMy00OTY3LWEyMmMtNTExZTQzNmY0OTdiL2ltYWdlXkEyXkFqcGdeQXVyNjQwOTYyNTY@._V1_SY1000_SX750_AL_.jpg" `
    "Cartoon Network" $true

Add-CartoonRow 23 "Chowder" "C. H. Greenblatt" `
    "In Marzipan City, the young, excitable food-loving Chowder is the apprentice of Mung Daal, a very old chef who runs a catering company with his wife, Truffles and assistant, Shnitzel." `
    2007 2010 "Comedy, Slapstick" 3 93 `
    "https://m.media-amazon.com/images/M/MV5BYmFlMDZkMWItZjVkMi00MzhmLWFhMjUtYTRjNjE1ZjE5ZGMzXkEyXkFqcGdeQXVyNTAyODkwOQ@@._V1_SY1000_CR0,0,666,1000_AL_.jpg" `
    "Cartoon Network" $true

# ---------------------------------------------------------------------
# 3. Column width / layout tweaks.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 39.256
$ws.Columns.Item(5).ColumnWidth = 7.893229166666667

# ---------------------------------------------------------------------
# 4. Selection moves to F9, matching the saved workbook state.
# ---------------------------------------------------------------------
$ws.Range("F9").Select() | Out-Null
